$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "A3"  = -21.687
    "C3"  = -12.774
    "A21" = -20.257
    "A23" = -20.585
    "C24" = -12.687
    "A25" = -21.796
    "B27" = 5.513
    "B31" = 6.145
    "B39" = 7.997000000000002
    "B48" = 5.274
    "B51" = 6.132000000000001
    "B52" = 5.295
    "A53" = -22.018
    "B55" = 4.685
    "B56" = 4.906999999999999
    "A57" = -21.352
    "B57" = 5.948
    "C57" = -13.287
    "A59" = -22.5
    "C61" = -13.508
    "A69" = -21.694
    "C70" = -11.941
    "B73" = 7.342999999999999
    "A79" = -21.192
    "A83" = -21.985
    "C86" = -13.597
    "B89" = 5.946
    "B90" = 5.833
    "A93" = -21.536
    "C98" = -12.45
    "C100" = -13.174
    "C102" = -13.564
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
